# Fix hyperlink missing after save issue and update completion row element
# logic: append a new row below the existing data with a "GitHub" label
# that links out to the project's repository.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row (row 22) with the display text for the hyperlink.
$ws.Range("A22").Value = "GitHub"

# Attach the hyperlink to the cell we just populated.
$ws.Hyperlinks.Add($ws.Range("A22"), "https://github.com/xuri/excelize")

# Hyperlinks.Add auto-applies Excel's built-in "Hyperlink" style (blue /
# underlined); restore the default cell style so the cell keeps its plain
# formatting, matching the original edit.
$ws.Range("A22").Style = "Normal"
